# SCD0011-017 - Penyelia Mengakses Menu Report - Menu Product Holding Ratio - Report
# Update sheet "SCD0186" -> "SCD0011" and TC_ID "DGS-201" -> "SCD0011-017"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (was "SCD0186")
$ws.Name = "SCD0011"

# Update the TC_ID cell value (was "DGS-201")
$ws.Range("B2").Value = "SCD0011-017"

# Widen column B to fit the new, longer TC_ID text
$ws.Range("B1").ColumnWidth = 11.6667

# Move the active selection to B3
$ws.Range("B3").Select() | Out-Null
